# Append a new data row (row 12, "2021年") to Sheet1, following the same
# layout/pattern as the existing rows (row 2 .. row 11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 12
$sourceRow = 11

# Copy the formatting of the label cell in column A from the last existing
# row so the new year label keeps the same (bold / bordered / centered)
# style used by every other row in column A.
$ws.Range("A" + $sourceRow).Copy() | Out-Null
$ws.Range("A" + $newRow).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Year label for the new row.
$ws.Range("A" + $newRow).Value = "2021年"

# Numeric data for the new row, keyed by column letter.
$rowValues = @{
    "B"  = 1
    "C"  = 308
    "F"  = 34
    "H"  = 84
    "L"  = 98
    "O"  = 21
    "Q"  = 5
    "U"  = 527
    "V"  = 42
    "X"  = 5773
    "Y"  = 1
    "AA" = 301
    "AB" = 7425
    "AD" = 107
    "AH" = 4
    "AJ" = 35
    "AL" = 12
    "AM" = 11
    "AO" = 24
    "AV" = 37
}

foreach ($col in $rowValues.Keys) {
    $ws.Range($col + $newRow).Value = $rowValues[$col]
}

Write-Host "Row 12 (2021年) added."
